$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '98.418.58'
Set-TextValue 'E2' '  +4.66%  '
Set-TextValue 'D3' '3.361.70'
Set-TextValue 'E3' '  +9.22%  '
Set-TextValue 'E4' '  +0.06%  '
Set-TextValue 'D5' '256.40'
Set-TextValue 'E5' '  +8.35%  '
Set-TextValue 'D6' '623.56'
Set-TextValue 'E6' '  +2.48%  '
Set-TextValue 'D7' '1.22'
Set-TextValue 'E7' '  +9.58%  '
Set-TextValue 'D8' '0.387'
Set-TextValue 'E8' '  +1.71%  '
Set-TextValue 'E9' '  +0.03%  '
Set-TextValue 'D10' '3.357.76'
Set-TextValue 'E10' '  +9.24%  '
Set-TextValue 'D11' '0.816'
Set-TextValue 'E11' '  +1.24%  '
Set-TextValue 'E12' '  +1.57%  '
Set-TextValue 'D13' '98.052.89'
Set-TextValue 'E13' '  +4.67%  '
Set-TextValue 'D14' '35.77'
Set-TextValue 'E14' '  +5.70%  '
Set-TextValue 'D15' '0.0000248'
Set-TextValue 'E15' '  +2.83%  '
Set-TextValue 'D16' '3.997.07'
Set-TextValue 'E16' '  +9.76%  '
Set-TextValue 'D17' '5.50'
Set-TextValue 'E17' '  +3.65%  '
Set-TextValue 'D18' '3.364.72'
Set-TextValue 'E18' '  +9.37%  '
Set-TextValue 'D19' '3.66'
Set-TextValue 'E19' '  +2.74%  '
Set-TextValue 'D20' '15.01'
Set-TextValue 'E20' '  +4.90%  '
Set-TextValue 'D21' '486.89'
Set-TextValue 'E21' '  +10.29%  '
Set-TextValue 'D22' '5.92'
Set-TextValue 'E22' '  +3.01%  '
Set-TextValue 'E23' '  +9.74%  '
Set-TextValue 'D24' '9.28'
Set-TextValue 'E24' '  +4.76%  '
Set-TextValue 'E25' '  +3.51%  '
Set-TextValue 'D26' '88.07'
Set-TextValue 'E26' '  +4.23%  '
Set-TextValue 'D27' '12.07'
Set-TextValue 'D28' '3.543.86'
Set-TextValue 'D29' '0.999'
Set-TextValue 'E29' '  +0.01%  '
Set-TextValue 'D30' '0.258'
Set-TextValue 'E30' '  +4.08%  '
Set-TextValue 'D31' '0.187'
Set-TextValue 'E31' '  +4.36%  '
Set-TextValue 'D32' '0.126'
Set-TextValue 'E32' '  +2.42%  '
Set-TextValue 'E33' '  +0.13%  '
Set-TextValue 'E34' '  +3.89%  '
Set-TextValue 'D35' '27.44'
Set-TextValue 'E35' '  +7.63%  '
Set-TextValue 'D36' '522.11'
Set-TextValue 'E36' '  +7.54%  '
Set-TextValue 'D37' '0.153'
Set-TextValue 'E37' '  -0.15%  '
Set-TextValue 'D38' '7.37'
Set-TextValue 'E38' '  -1.40%  '
Set-TextValue 'D40' '24.80'
Set-TextValue 'D41' '0.451'
Set-TextValue 'E41' '  +3.40%  '
Set-TextValue 'D42' '3.76'
Set-TextValue 'E42' '  -2.19%  '
Set-TextValue 'D43' '1.27'
Set-TextValue 'E43' '  +2.33%  '
Set-TextValue 'D44' '3.26'
Set-TextValue 'E44' '  +5.64%  '
Set-TextValue 'D45' '0.782'
Set-TextValue 'E45' '  +15.70%  '
Set-TextValue 'E46' '  -0.04%  '
Set-TextValue 'D47' '160.39'
Set-TextValue 'E47' '  -0.48%  '
Set-TextValue 'E48' '  +5.64%  '
Set-TextValue 'D49' '4.55'
Set-TextValue 'E49' '  +6.62%  '
Set-TextValue 'D50' '45.43'
Set-TextValue 'E51' '  +6.16%  '
